# Updates the crypto price/volume snapshot (GitHub Actions refresh).
# Every Coin/Link/Price/Volume cell in the sheet is stored as literal
# text (prices use '.' as a thousands separator, e.g. '26.079.86', and
# volumes are padded strings like '  +1.23%  '). Plain decimal prices
# (e.g. '63.20', '10.00', '0.416') read back as numbers through
# Range.Value's automatic type detection, which would silently drop
# exact trailing zeros, so those are written with a leading apostrophe
# to force literal text, matching the source data exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.079.86'
$ws.Range('E2').Value = '  +1.23%  '
$ws.Range('D3').Value = '1.638.81'
$ws.Range('E3').Value = '  +0.31%  '
$ws.Range('E4').Value = '  -0.30%  '
$ws.Range('D5').Value = '''216.22'
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('E6').Value = '  +0.50%  '
$ws.Range('D7').Value = '''0.999'
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('D8').Value = '''0.257'
$ws.Range('E8').Value = '  +0.27%  '
$ws.Range('D9').Value = '''0.0634'
$ws.Range('E9').Value = '  +0.13%  '
$ws.Range('E10').Value = '  +0.63%  '
$ws.Range('E11').Value = '  -0.16%  '
$ws.Range('D12').Value = '1.867.59'
$ws.Range('E12').Value = '  +0.42%  '
$ws.Range('D13').Value = '''4.27'
$ws.Range('E13').Value = '  +0.28%  '
$ws.Range('D14').Value = '1.600.13'
$ws.Range('E14').Value = '  -2.22%  '
$ws.Range('E15').Value = '  -2.85%  '
$ws.Range('D16').Value = '0.0₃0761'
$ws.Range('E16').Value = '  -0.30%  '
$ws.Range('D17').Value = '''63.20'
$ws.Range('E17').Value = '  +0.15%  '
$ws.Range('D18').Value = '26.091.16'
$ws.Range('E18').Value = '  +1.19%  '
$ws.Range('E19').Value = '  -0.26%  '
$ws.Range('D20').Value = '''194.86'
$ws.Range('E20').Value = '  +1.32%  '
$ws.Range('E21').Value = '  -0.87%  '
$ws.Range('D22').Value = '''10.00'
$ws.Range('E22').Value = '  +0.22%  '
$ws.Range('D23').Value = '''6.26'
$ws.Range('E23').Value = '  -0.37%  '
$ws.Range('D24').Value = '''1.78'
$ws.Range('E24').Value = '  -1.31%  '
$ws.Range('D25').Value = '''0.999'
$ws.Range('E25').Value = '  -0.27%  '
$ws.Range('D26').Value = '''142.54'
$ws.Range('E26').Value = '  +0.16%  '
$ws.Range('E27').Value = '  +0.76%  '
$ws.Range('E28').Value = '  +0.32%  '
$ws.Range('D29').Value = '''15.57'
$ws.Range('E29').Value = '  +0.57%  '
$ws.Range('E30').Value = '  +0.10%  '
$ws.Range('E31').Value = '  +2.09%  '
$ws.Range('E32').Value = '  +0.20%  '
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('D34').Value = '''1.58'
$ws.Range('E34').Value = '  +0.84%  '
$ws.Range('E35').Value = '  +1.40%  '
$ws.Range('D36').Value = '''0.910'
$ws.Range('E36').Value = '  +0.57%  '
$ws.Range('D37').Value = '1.134.22'
$ws.Range('E37').Value = '  +0.21%  '
$ws.Range('D38').Value = '''0.551'
$ws.Range('E38').Value = '  +1.04%  '
$ws.Range('E39').Value = '  -1.20%  '
$ws.Range('E40').Value = '  +1.14%  '
$ws.Range('D41').Value = '''0.998'
$ws.Range('E41').Value = '  -0.35%  '
$ws.Range('D42').Value = '''99.83'
$ws.Range('E42').Value = '  -0.86%  '
$ws.Range('E43').Value = '  -1.54%  '
$ws.Range('E44').Value = '  -0.20%  '
$ws.Range('D45').Value = '1.776.74'
$ws.Range('E45').Value = '  +0.44%  '
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('D47').Value = '''56.62'
$ws.Range('E47').Value = '  +2.15%  '
$ws.Range('E48').Value = '  +2.15%  '
$ws.Range('E49').Value = '  +3.31%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').Value = '''0.416'
$ws.Range('E50').Value = '  -0.16%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '''7.66'
$ws.Range('E51').Value = '  +2.48%  '
